$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.071549776775584
$ws.Range("D2").Value = 1.074613030876448
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.083501457501438
$ws.Range("I2").Value = 1.057250105275551
$ws.Range("J2").Value = 1.076473808171504
$ws.Range("K2").Value = 1.077301780854576
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.086166909069458
$ws.Range("N2").Value = 1.078002524452379
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.072786702988307
$ws.Range("D3").Value = 1.07561883176587
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.084638274852635
$ws.Range("I3").Value = 1.057654343875325
$ws.Range("J3").Value = 1.07736734567015
$ws.Range("K3").Value = 1.07812427690539
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.087121789610901
$ws.Range("N3").Value = 1.078897330876763
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.073586577157994
$ws.Range("D4").Value = 1.07626918190617
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.085373697675744
$ws.Range("I4").Value = 1.057914461699655
$ws.Range("J4").Value = 1.077944466793571
$ws.Range("K4").Value = 1.078655401353794
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.087738895706541
$ws.Range("N4").Value = 1.079475271578376
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07392272692633
$ws.Range("D5").Value = 1.076542477744614
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.085682829196891
$ws.Range("I5").Value = 1.058023468432619
$ws.Range("J5").Value = 1.078186836839003
$ws.Range("K5").Value = 1.078878427445401
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.08799814519528
$ws.Range("N5").Value = 1.079717985817073
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.073979161138535
$ws.Range("D6").Value = 1.076588358823837
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.085734731421628
$ws.Range("I6").Value = 1.058041750826837
$ws.Range("J6").Value = 1.078227517120324
$ws.Range("K6").Value = 1.07891585937764
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.088041663696569
$ws.Range("N6").Value = 1.07975872386906
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.073591069265784
$ws.Range("D7").Value = 1.076272834134578
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.085377828462567
$ws.Range("I7").Value = 1.057915919614757
$ws.Range("J7").Value = 1.077947706341062
$ws.Range("K7").Value = 1.078658382454021
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.087742360524022
$ws.Range("N7").Value = 1.079478515726396
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.071967906665483
$ws.Range("D8").Value = 1.074953044113938
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.083885687140306
$ws.Range("I8").Value = 1.057387020590776
$ws.Range("J8").Value = 1.07677600317532
$ws.Range("K8").Value = 1.077579972754459
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.086489774721989
$ws.Range("N8").Value = 1.078305148607763
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.069103754651956
$ws.Range("D9").Value = 1.072623748045478
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.081254946254327
$ws.Range("I9").Value = 1.056443883474534
$ws.Range("J9").Value = 1.074703148103896
$ws.Range("K9").Value = 1.075671307361818
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.084276635445957
$ws.Range("N9").Value = 1.076229349844378
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.067191528128507
$ws.Range("D10").Value = 1.071068344574764
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.07950008668285
$ws.Range("I10").Value = 1.055807582885756
$ws.Range("J10").Value = 1.073315660798892
$ws.Range("K10").Value = 1.074393155708833
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.082797139848423
$ws.Range("N10").Value = 1.074839892148253
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.0663628170482
$ws.Range("D11").Value = 1.070394216653963
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.078739946722497
$ws.Range("I11").Value = 1.055530257986799
$ws.Range("J11").Value = 1.072713517264057
$ws.Range("K11").Value = 1.07383832937649
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.08215551654785
$ws.Range("N11").Value = 1.074236893500532
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.066054888026658
$ws.Range("D12").Value = 1.070143719468411
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.078457553741581
$ws.Range("I12").Value = 1.055426975365256
$ws.Range("J12").Value = 1.072489649190558
$ws.Range("K12").Value = 1.073632033402876
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.081917038218707
$ws.Range("N12").Value = 1.074012707508691
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.066120944820091
$ws.Range("D13").Value = 1.070197456339642
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.07851812997901
$ws.Range("I13").Value = 1.055449142153978
$ws.Range("J13").Value = 1.07253767896772
$ws.Range("K13").Value = 1.073676294064271
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.081968199486504
$ws.Range("N13").Value = 1.074060805493645
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.06633736577668
$ws.Range("D14").Value = 1.070373512448025
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.078716604922313
$ws.Range("I14").Value = 1.055521726163607
$ws.Range("J14").Value = 1.072695016455799
$ws.Range("K14").Value = 1.073821281155699
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.082135806943073
$ws.Range("N14").Value = 1.074218366419004
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.066470695360075
$ws.Range("D15").Value = 1.07048197364656
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.078838886108192
$ws.Range("I15").Value = 1.055566411520768
$ws.Range("J15").Value = 1.072791930030866
$ws.Range("K15").Value = 1.073910584772675
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.082239055382331
$ws.Range("N15").Value = 1.074315417622465
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.067246511857068
$ws.Range("D16").Value = 1.071113070848035
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.079550528724595
$ws.Range("I16").Value = 1.055825949963614
$ws.Range("J16").Value = 1.073355594420776
$ws.Range("K16").Value = 1.074429948521984
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.082839701217869
$ws.Range("N16").Value = 1.07487988248046
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.067732970349151
$ws.Range("D17").Value = 1.071508772117877
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.079996848546906
$ws.Range("I17").Value = 1.055988268278223
$ws.Range("J17").Value = 1.073708802738863
$ws.Range("K17").Value = 1.074755361485738
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.083216203613989
$ws.Range("N17").Value = 1.075233592394874
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.068016645599695
$ws.Range("D18").Value = 1.071739517502985
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.080257153027346
$ws.Range("I18").Value = 1.056082771953352
$ws.Range("J18").Value = 1.073914692825141
$ws.Range("K18").Value = 1.074945036499127
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.083435715509237
$ws.Range("N18").Value = 1.07543977486869
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.068113360126674
$ws.Range("D19").Value = 1.071818185510773
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.080345905722553
$ws.Range("I19").Value = 1.056114965805932
$ws.Range("J19").Value = 1.073984873973522
$ws.Range("K19").Value = 1.075009688292339
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.083510547235448
$ws.Range("N19").Value = 1.075510055682351
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.067680784988225
$ws.Range("D20").Value = 1.071466323384684
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.079948965391641
$ws.Range("I20").Value = 1.055970871049232
$ws.Range("J20").Value = 1.073670920327005
$ws.Range("K20").Value = 1.074720461510257
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.083175818370887
$ws.Range("N20").Value = 1.075195656185645
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.066273638227467
$ws.Range("D21").Value = 1.070321670977757
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.078658160205437
$ws.Range("I21").Value = 1.055500359493061
$ws.Range("J21").Value = 1.072648690166255
$ws.Range("K21").Value = 1.073778591847799
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.082086454904095
$ws.Range("N21").Value = 1.074171974340815
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.065388278070815
$ws.Range("D22").Value = 1.069601425947123
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.077846329097463
$ws.Range("I22").Value = 1.055202957121993
$ws.Range("J22").Value = 1.072004785549655
$ws.Range("K22").Value = 1.073185192278749
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.081400655069674
$ws.Range("N22").Value = 1.073527155305803
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.06585768507169
$ws.Range("D23").Value = 1.069983294803458
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.078276720544814
$ws.Range("I23").Value = 1.05536076518875
$ws.Range("J23").Value = 1.072346244834908
$ws.Range("K23").Value = 1.073499879763228
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.081764293978745
$ws.Range("N23").Value = 1.073869099502407
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.067704365515865
$ws.Range("D24").Value = 1.071485504326939
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.07997060180979
$ws.Range("I24").Value = 1.055978732645346
$ws.Range("J24").Value = 1.073688038161068
$ws.Range("K24").Value = 1.074736231717478
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.083194067019258
$ws.Range("N24").Value = 1.075212798328996
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.069844687800476
$ws.Range("D25").Value = 1.073226368583161
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.08193523094691
$ws.Range("I25").Value = 1.056689033045313
$ws.Range("J25").Value = 1.075240008069137
$ws.Range("K25").Value = 1.076165743259887
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.084849495850852
$ws.Range("N25").Value = 1.076766972212349
